{"js": "const replacements = [\n  [\"2024-12-20 Friday\", \"2024-12-21 Saturday\"],\n  [\"84\u00f78=\", \"61\u00f74=\"],\n  [\"89\u00f74=\", \"68\u00f76=\"],\n  [\"77\u00f77=\", \"17\u00f76=\"],\n  [\"78\u00f78=\", \"94\u00f79=\"],\n  [\"85\u00f78=\", \"94\u00f72=\"],\n  [\"95\u00f76=\", \"57\u00f77=\"],\n  [\"48\u00f73=\", \"87\u00f79=\"],\n  [\"81\u00f79=\", \"31\u00f76=\"],\n  [\"97\u00f72=\", \"10\u00f78=\"],\n  [\"89\u00f73=\", \"66\u00f79=\"],\n  [\"75\u00f72=\", \"93\u00f76=\"],\n  [\"11\u00f76=\", \"72\u00f76=\"],\n  [\"23\u00f76=\", \"64\u00f77=\"],\n  [\"27\u00f78=\", \"95\u00f74=\"],\n  [\"63\u00f79=\", \"39\u00f72=\"],\n  [\"91\u00f75=\", \"69\u00f77=\"],\n  [\"67\u00f76=\", \"42\u00f76=\"],\n  [\"56\u00f74=\", \"65\u00f77=\"],\n  [\"79\u00f74=\", \"98\u00f75=\"],\n  [\"71\u00f74=\", \"59\u00f77=\"],\n  [\"90\u00f79=\", \"79\u00f75=\"],\n  [\"95\u00f79=\", \"60\u00f76=\"],\n  [\"36\u00f77=\", \"53\u00f72=\"],\n  [\"13\u00f78=\", \"79\u00f72=\"],\n  [\"19\u00f74=\", \"91\u00f73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-12-20 Friday\", \"2024-12-21 Saturday\"),\n  @(\"84\u00f78=\", \"61\u00f74=\"),\n  @(\"89\u00f74=\", \"68\u00f76=\"),\n  @(\"77\u00f77=\", \"17\u00f76=\"),\n  @(\"78\u00f78=\", \"94\u00f79=\"),\n  @(\"85\u00f78=\", \"94\u00f72=\"),\n  @(\"95\u00f76=\", \"57\u00f77=\"),\n  @(\"48\u00f73=\", \"87\u00f79=\"),\n  @(\"81\u00f79=\", \"31\u00f76=\"),\n  @(\"97\u00f72=\", \"10\u00f78=\"),\n  @(\"89\u00f73=\", \"66\u00f79=\"),\n  @(\"75\u00f72=\", \"93\u00f76=\"),\n  @(\"11\u00f76=\", \"72\u00f76=\"),\n  @(\"23\u00f76=\", \"64\u00f77=\"),\n  @(\"27\u00f78=\", \"95\u00f74=\"),\n  @(\"63\u00f79=\", \"39\u00f72=\"),\n  @(\"91\u00f75=\", \"69\u00f77=\"),\n  @(\"67\u00f76=\", \"42\u00f76=\"),\n  @(\"56\u00f74=\", \"65\u00f77=\"),\n  @(\"79\u00f74=\", \"98\u00f75=\"),\n  @(\"71\u00f74=\", \"59\u00f77=\"),\n  @(\"90\u00f79=\", \"79\u00f75=\"),\n  @(\"95\u00f79=\", \"60\u00f76=\"),\n  @(\"36\u00f77=\", \"53\u00f72=\"),\n  @(\"13\u00f78=\", \"79\u00f72=\"),\n  @(\"19\u00f74=\", \"91\u00f73=\"),\n)\n\nforeach ($p in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $p[0]\n  $find.Replacement.Text = $p[1]\n  $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
